# Dashboard de Pagamentos - tela de conciliacao bancaria
# The "company" detail row (row 2, with the "Empresa Modelo" label and its
# Pagamentos/Debitos/Total figures) is replaced by the previously-hidden
# footer-style row (old row 3), now showing zeroed-out totals, and the old
# footer row is removed from the sheet entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the footer row's formatting (style ids 1/2/1 pattern) onto row 2,
# replacing the numeric "Empresa Modelo" totals formatting.
$ws.Range("A3:E3").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)

# Row 2 becomes the blank-label / all-zero text row.
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "0,00"
$ws.Range("C2").Value = "0,00"
$ws.Range("D2").Value = "0,00"
$ws.Range("E2").Value = ""

# The old footer row (row 3) is no longer needed - remove it so the used
# range collapses back down to A1:E2.
$ws.Rows(3).Delete()
